# Fixed #295 Add the version of M2Doc in the template custom properties.
#
# The canonical-OOXML diff for this particular resource
# (extraSpaceInComment-template.docx) is entirely cosmetic: every
# hunk only re-orders XML attributes that some parts of the same
# commit's save pipeline re-serialized (e.g. w:pgSz w:w/w:h ->
# w:h/w:w, w:latentStyles/w:lsdException/w:style attribute order,
# the root <w:document> namespace-declaration order, ...). Tag
# names, attribute names, attribute values, element structure and
# visible content are all unchanged - every "-" line and its
# matching "+" line carry the exact same set of attribute
# name/value pairs, just written in a different order.
#
# Word's automation object model intentionally does not expose raw
# attribute-serialization order: PageSetup, Styles, Sections, etc.
# let you read/write the *values* (margins, page size, font theme,
# priorities, ...), and the document's XML is re-emitted with those
# values unchanged, but there is no supported COM property/method
# that reorders the attributes Word already wrote. So the
# content-faithful way to "apply" this diff through COM automation
# is to touch the document without altering any value - which is
# exactly what happens below: every figure is read back and written
# back unchanged, so the document keeps the same text, the same
# page setup and the same styles it started with.

$d = $word.ActiveDocument

# Page setup (w:sectPr / w:pgSz / w:pgMar in word/document.xml) -
# read back and reassign the same values; no margin/size actually
# changes.
$section = $d.Sections(1)
$pageSetup = $section.PageSetup
$pageSetup.PageWidth = $pageSetup.PageWidth
$pageSetup.PageHeight = $pageSetup.PageHeight
$pageSetup.TopMargin = $pageSetup.TopMargin
$pageSetup.BottomMargin = $pageSetup.BottomMargin
$pageSetup.LeftMargin = $pageSetup.LeftMargin
$pageSetup.RightMargin = $pageSetup.RightMargin
$pageSetup.HeaderDistance = $pageSetup.HeaderDistance
$pageSetup.FooterDistance = $pageSetup.FooterDistance
$pageSetup.Gutter = $pageSetup.Gutter

# Styles (word/styles.xml docDefaults/latentStyles/w:style entries) -
# same idea: re-read and reassign so every style keeps its existing
# priority/visibility/name, nothing is added, removed or renamed.
foreach ($style in $d.Styles) {
    $null = $style.NameLocal
    $null = $style.Priority
}
